$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B213").Value = "Story points "
$ws.Range("C213").Value = "Min."
$ws.Range("D213").Value = "Estimado"
$ws.Range("E213").Value = "Total"
$ws.Range("F213").Value = "Real"
$ws.Range("G213").Value = "Total3"
$lo = $ws.ListObjects.Add(1, $ws.Range("B213:G221"), 0, 1)
$lo.ShowTotals = $true
Write-Output ("post showtotals Range: " + $lo.Range.Address())
Write-Output ("post showtotals AutoFilter.Range: " + $lo.AutoFilter.Range.Address())

# try resetting autofilter off then on
$lo.ShowAutoFilter = $false
Write-Output ("autofilter off -> ShowAutoFilter=" + $lo.ShowAutoFilter)
$lo.ShowAutoFilter = $true
Write-Output ("autofilter on again -> Range: " + $lo.AutoFilter.Range.Address())
